$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "The Mirthless" boss to "The Acharos" (track name + description)
$ws.Range("A6").Value = "The Acharos"
$ws.Range("B6").Value = "Plays when the Acharos fights the player. The Acharos often haunts players even as they only begin to become insane, but at low sanity levels, it finally confronts the player."

# Update the sheet view: scroll so row 2 is at top, and select B6
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("B6").Select()
